$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.218.08"
$ws.Range("E2").Value = "  -2.12%  "

$ws.Range("D3").Value = "2.899.03"
$ws.Range("E3").Value = "  -3.20%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").Value = "529.14"
$ws.Range("E5").Value = "  -1.61%  "

$ws.Range("D6").Value = "130.79"
$ws.Range("E6").Value = "  -1.55%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "2.896.74"
$ws.Range("E8").Value = "  -3.11%  "

$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  -3.82%  "

$ws.Range("D10").Value = "6.28"
$ws.Range("E10").Value = "  +3.08%  "

$ws.Range("D11").Value = "0.141"
$ws.Range("E11").Value = "  -3.87%  "

$ws.Range("D12").Value = "0.428"
$ws.Range("E12").Value = "  -4.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000212"
$ws.Range("E13").Value = "  -4.40%  "

$ws.Range("D14").Value = "32.38"
$ws.Range("E14").Value = "  -3.99%  "

$ws.Range("D15").Value = "3.415.38"
$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("D16").Value = "60.267.66"
$ws.Range("E16").Value = "  -2.22%  "

$ws.Range("D17").Value = "0.106"
$ws.Range("E17").Value = "  -3.62%  "

$ws.Range("D18").Value = "2.915.91"
$ws.Range("E18").Value = "  -2.83%  "

$ws.Range("D19").Value = "6.37"
$ws.Range("E19").Value = "  -3.58%  "

$ws.Range("D20").Value = "451.95"
$ws.Range("E20").Value = "  -3.67%  "

$ws.Range("D21").Value = "12.91"
$ws.Range("E21").Value = "  -1.72%  "

$ws.Range("D22").Value = "0.627"
$ws.Range("E22").Value = "  -6.08%  "

$ws.Range("D23").Value = "6.81"
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("D24").Value = "77.06"
$ws.Range("E24").Value = "  -3.85%  "

$ws.Range("D25").Value = "11.94"
$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "2.62"
$ws.Range("E27").Value = "  -2.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.35%  "

$ws.Range("D29").Value = "7.24"
$ws.Range("E29").Value = "  -5.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.90"
$ws.Range("E30").Value = "  +0.90%  "

$ws.Range("D31").Value = "24.45"
$ws.Range("E31").Value = "  -4.08%  "

$ws.Range("D32").Value = "1.09"
$ws.Range("E32").Value = "  -4.95%  "

$ws.Range("D33").Value = "2.24"
$ws.Range("E33").Value = "  -2.11%  "

$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.00"
$ws.Range("E35").Value = "  -4.53%  "

$ws.Range("D36").Value = "5.62"
$ws.Range("E36").Value = "  -4.36%  "

$ws.Range("D37").Value = "432.57"
$ws.Range("E37").Value = "  -4.96%  "

$ws.Range("D38").Value = "0.0772"
$ws.Range("E38").Value = "  -2.01%  "

$ws.Range("D39").Value = "0.0373"
$ws.Range("E39").Value = "  -2.02%  "

$ws.Range("D40").Value = "2.857.87"
$ws.Range("E40").Value = "  -10.14%  "

$ws.Range("D41").Value = "0.111"
$ws.Range("E41").Value = "  -6.90%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.70"
$ws.Range("E42").Value = "  -4.58%  "

$ws.Range("E43").Value = "  -0.08%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.30"
$ws.Range("E44").Value = "  -3.11%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "25.35"
$ws.Range("E45").Value = "  -0.68%  "

$ws.Range("D46").Value = "0.237"
$ws.Range("E46").Value = "  -1.86%  "

$ws.Range("D47").Value = "0.105"
$ws.Range("E47").Value = "  -2.06%  "

$ws.Range("D48").Value = "1.88"
$ws.Range("E48").Value = "  -4.40%  "

$ws.Range("D49").Value = "111.39"
$ws.Range("E49").Value = "  -5.97%  "

$ws.Range("D50").Value = "0.0₃0466"
$ws.Range("E50").Value = "  -4.27%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.20"
$ws.Range("E51").Value = "  -3.88%  "
